# CertificateNR35.pptx — "remove RG to user and models certificate"
#
# Slide 1, shape "Rectangle 5" currently reads (paragraph 1):
#   Certificamos que {{NOME}}, portador do RG nº {RG}} e CPF nº {{CPF}}, concluiu ...
# and must become:
#   Certificamos que {{NOME}}, portador do CPF nº {{CPF}}, concluiu ...
# i.e. the "RG nº {RG}} e " part is dropped, leaving only the CPF reference.
#
# The shape uses <a:spAutoFit/>, so once the paragraph text is edited the
# runtime recomputes the shape's rendered height on its own (matching the
# <a:ext cy="…"/> shrink in the diff) — no manual resize is required.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Rectangle 5")
$tr = $sh.TextFrame.TextRange

# Original (pre-edit) character layout of paragraph 1, all 1-based:
#   26-27  ", "
#   28-45  "portador do RG nº "
#   46-51  "{RG}} "
#   52-60  "e CPF nº "
#   61-69  "{{CPF}}, "
#
# Rewrite right-to-left so earlier offsets stay valid while later ones shift.
$tr.Characters(61, 9).Text  = "nº {{CPF}}, "
$tr.Characters(46, 15).Text = "CPF "
$tr.Characters(28, 18).Text = "do "
$tr.Characters(26, 2).Text  = ", portador "
